$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 50149
$ws.Range("I12").Value = 50149
$ws.Range("K12").Value = 50149
$ws.Range("M12").Value = -49979

$ws.Range("H40").Value = 3745.7144
$ws.Range("I40").Value = 2215.5557
$ws.Range("K40").Value = 2215.5557
$ws.Range("M40").Value = -2040.5557

$ws.Range("H88").Value = 16781912
$ws.Range("J88").Value = 2756982.5
$ws.Range("L88").Value = 2756982.5
$ws.Range("N88").Value = -2757794.5

$ws.Range("H91").Value = 16781912
$ws.Range("J91").Value = 2756982.5
$ws.Range("L91").Value = 2756982.5
$ws.Range("N91").Value = -2759790.5

$ws.Range("H131").Value = 3812.8
$ws.Range("I131").Value = 1516
$ws.Range("K131").Value = 4548
$ws.Range("M131").Value = 492

$ws.Range("H137").Value = 2477.3572
$ws.Range("I137").Value = 1618.4
$ws.Range("K137").Value = 4855.200000000001
$ws.Range("M137").Value = -2305.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 41669188
$ws.Range("I61").Value = 45456204
$ws.Range("K61").Value = 45456204
$ws.Range("M61").Value = -45455992

$ws.Range("H76").Value = 52499.5
$ws.Range("I76").Value = 50000
$ws.Range("K76").Value = 50000
$ws.Range("M76").Value = -49662

$ws.Range("H79").Value = 52499.5
$ws.Range("I79").Value = 50000
$ws.Range("K79").Value = 50000
$ws.Range("M79").Value = -48830

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H98").Value = 42660
$ws.Range("J98").Value = 42660
$ws.Range("L98").Value = 42660
$ws.Range("N98").Value = -48650

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 41669188
$ws.Range("I136").Value = 45456204
$ws.Range("K136").Value = 136368612
$ws.Range("M136").Value = -136366062

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 11965
$ws.Range("I82").Value = 11965
$ws.Range("K82").Value = 11965
$ws.Range("M82").Value = -11582

$ws.Range("H85").Value = 11965
$ws.Range("I85").Value = 11965
$ws.Range("K85").Value = 11965
$ws.Range("M85").Value = -10639

$ws.Range("H94").Value = 3205.1
$ws.Range("I94").Value = 3406.1333
$ws.Range("J94").Value = 2602
$ws.Range("K94").Value = 3406.1333
$ws.Range("L94").Value = 2602
$ws.Range("M94").Value = -2955.1333
$ws.Range("N94").Value = -3504

$ws.Range("H103").Value = 22564
$ws.Range("J103").Value = 22564
$ws.Range("L103").Value = 22564
$ws.Range("N103").Value = -24908

$ws.Range("H134").Value = 16455368
$ws.Range("I134").Value = 18218050
$ws.Range("K134").Value = 54654150
$ws.Range("M134").Value = -54651615

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 25999
$ws.Range("J9").Value = 25999
$ws.Range("L9").Value = 25999
$ws.Range("N9").Value = -26335

$ws.Range("H22").Value = 7044.2
$ws.Range("I22").Value = 7982.385
$ws.Range("K22").Value = 7982.385
$ws.Range("M22").Value = -7632.385

$ws.Range("H23").Value = 5000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 5000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 5000
$ws.Range("N23").Value = -5480
$ws.Range("M23").ClearContents()

$ws.Range("H27").Value = 5000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 5000
$ws.Range("N27").Value = -5384
$ws.Range("M27").ClearContents()

$ws.Range("H31").Value = 3178.375
$ws.Range("I31").Value = 3187.353
$ws.Range("K31").Value = 3187.353
$ws.Range("M31").Value = -2892.353

$ws.Range("H34").Value = 3178.375
$ws.Range("I34").Value = 3187.353
$ws.Range("K34").Value = 3187.353
$ws.Range("M34").Value = -2985.353

$ws.Range("H134").Value = 83335000
$ws.Range("I134").Value = 83335000
$ws.Range("K134").Value = 250005000
$ws.Range("M134").Value = -250002465

$ws.Range("H141").Value = 84963.6
$ws.Range("J141").Value = 93704.5
$ws.Range("L141").Value = 93704.5
$ws.Range("N141").Value = -104064.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 557191.9
$ws.Range("I4").Value = 682690.06
$ws.Range("K4").Value = 2048070.18
$ws.Range("M4").Value = -2047958.18

$ws.Range("H7").Value = 1669183.4
$ws.Range("J7").Value = 2517
$ws.Range("L7").Value = 7551
$ws.Range("N7").Value = -7775

$ws.Range("H11").Value = 118027.42
$ws.Range("I11").Value = 118748.52
$ws.Range("K11").Value = 356245.56
$ws.Range("M11").Value = -356105.56

$ws.Range("H40").Value = 24.666666
$ws.Range("I40").Value = 24.666666
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 98.666664
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -29.666664
$ws.Range("N40").ClearContents()

$ws.Range("H94").Value = 13613.333
$ws.Range("I94").Value = 3500
$ws.Range("J94").Value = 18670
$ws.Range("K94").Value = 10500
$ws.Range("L94").Value = 56010
$ws.Range("M94").Value = -9824
$ws.Range("N94").Value = -57362

$ws.Range("H109").Value = 1780.875
$ws.Range("I109").Value = 1606.7142
$ws.Range("K109").Value = 4820.142599999999
$ws.Range("M109").Value = -3780.142599999999

$ws.Range("H114").Value = 202002
$ws.Range("J114").Value = 4610.5
$ws.Range("L114").Value = 13831.5
$ws.Range("N114").Value = -20339.5

$ws.Range("H134").Value = 498.85715
$ws.Range("I134").Value = 498.85715
$ws.Range("K134").Value = 1496.57145
$ws.Range("M134").Value = 3573.42855

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 10000
$ws.Range("I53").Value = 10000
$ws.Range("K53").Value = 10000
$ws.Range("M53").Value = -9369

$ws.Range("H98").Value = 128194
$ws.Range("J98").Value = 128194
$ws.Range("L98").Value = 128194
$ws.Range("N98").Value = -134184

$ws.Range("H132").Value = 8335712
$ws.Range("I132").Value = 13891248
$ws.Range("K132").Value = 41673744
$ws.Range("M132").Value = -41671214

$ws.Range("H136").Value = 55000
$ws.Range("J136").Value = 55000
$ws.Range("L136").Value = 165000
$ws.Range("N136").Value = -170100

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1247.0625
$ws.Range("I16").Value = 1304.8667
$ws.Range("J16").Value = 380
$ws.Range("K16").Value = 1304.8667
$ws.Range("L16").Value = 380
$ws.Range("M16").Value = -1134.8667
$ws.Range("N16").Value = -720

$ws.Range("I46").Value = 1249.4286
$ws.Range("J46").Value = 979.3333
$ws.Range("K46").Value = 1249.4286
$ws.Range("L46").Value = 979.3333
$ws.Range("M46").Value = -1061.4286
$ws.Range("N46").Value = -1355.3333

$ws.Range("H68").Value = 2824.5
$ws.Range("I68").Value = 2824.5
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2824.5
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2075.5
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 2824.5
$ws.Range("I71").Value = 2824.5
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 14122.5
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -10378.5
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 2457.75
$ws.Range("I82").Value = 2340
$ws.Range("J82").Value = 2541.8572
$ws.Range("K82").Value = 2340
$ws.Range("L82").Value = 2541.8572
$ws.Range("M82").Value = -1979
$ws.Range("N82").Value = -3263.8572

$ws.Range("H85").Value = 2457.75
$ws.Range("I85").Value = 2340
$ws.Range("J85").Value = 2541.8572
$ws.Range("K85").Value = 2340
$ws.Range("L85").Value = 2541.8572
$ws.Range("M85").Value = -1092
$ws.Range("N85").Value = -5037.8572

$ws.Range("H122").Value = 10689
$ws.Range("I122").Value = 3226.1428
$ws.Range("K122").Value = 9678.4284
$ws.Range("M122").Value = -7228.428400000001

$ws.Range("H132").Value = 19208938
$ws.Range("I132").Value = 19208938
$ws.Range("K132").Value = 57626814
$ws.Range("M132").Value = -57624284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2700.6
$ws.Range("I122").Value = 2700.6
$ws.Range("K122").Value = 8101.799999999999
$ws.Range("M122").Value = -5651.799999999999

$ws.Range("H132").Value = 17864294
$ws.Range("I132").Value = 27779696
$ws.Range("K132").Value = 83339088
$ws.Range("M132").Value = -83336558
